$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "NA" under duplicate_image_filename (column E) for data rows 2 through 21
for ($r = 2; $r -le 21; $r++) {
    $ws.Range("E$r").Value = "NA"
}
